$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.9123456478118896
$ws.Range("D2").Value = 17.044
$ws.Range("C3").Value = 0.4934020042419434
$ws.Range("D3").Value = 66.15900000000001
$ws.Range("C4").Value = 0.5929369926452637
$ws.Range("D4").Value = 32.984
$ws.Range("C5").Value = 0.6129200458526611
$ws.Range("D5").Value = 7.308
$ws.Range("C6").Value = 0.2025918960571289
$ws.Range("D6").Value = 6.017
$ws.Range("C7").Value = 0.1944859027862549
$ws.Range("D7").Value = 6.017
$ws.Range("C8").Value = 0.2353639602661133
$ws.Range("D8").Value = 5.353
$ws.Range("C9").Value = 0.2401340007781982
$ws.Range("D9").Value = 8.99
$ws.Range("C10").Value = 0.3002710342407227
$ws.Range("D10").Value = 69.571
$ws.Range("C11").Value = 0.2825248241424561
$ws.Range("D11").Value = 73.881
$ws.Range("C12").Value = 1.456033945083618
$ws.Range("D12").Value = 69.503
$ws.Range("C13").Value = 0.6088719367980957
$ws.Range("D13").Value = 73.245
$ws.Range("C14").Value = 0.485663890838623
$ws.Range("D14").Value = 8.17
$ws.Range("C15").Value = 0.5502300262451172
$ws.Range("D15").Value = 8.17
$ws.Range("C16").Value = 0.5654866695404053
$ws.Range("D16").Value = 7.326
$ws.Range("C17").Value = 0.6181249618530273
$ws.Range("D17").Value = 7.326
$ws.Range("C18").Value = 0.2139163017272949
$ws.Range("D18").Value = 6.215
$ws.Range("C19").Value = 0.2209579944610596
$ws.Range("D19").Value = 6.215
$ws.Range("C20").Value = 0.2463538646697998
$ws.Range("D20").Value = 5.371
$ws.Range("C21").Value = 0.2633438110351562
$ws.Range("D21").Value = 5.371
$ws.Range("C22").Value = 0.2959749698638916
$ws.Range("D22").Value = 64.398
$ws.Range("C23").Value = 0.3065907955169678
$ws.Range("D23").Value = 522.806
$ws.Range("C24").Value = 0.6373012065887451
$ws.Range("D24").Value = 58.093
$ws.Range("C25").Value = 0.7557830810546875
$ws.Range("D25").Value = 522.668
